$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$words = @(
    'reiten',
    'kosten',
    'schreiten',
    'wählen',
    'fallen',
    'stehlen',
    'jagen',
    'lösen',
    'biegen',
    'weinen',
    'zeugen',
    'brauen',
    'knarren',
    'fischen',
    'drucken',
    'schlafen',
    'bluten',
    'trösten',
    'tanzen',
    'schulden',
    'sprechen',
    'irren',
    'filmen',
    'zögern',
    'mauern',
    'schultern',
    'sperren',
    'spinnen',
    'heißen',
    'fordern',
    'grüßen',
    'sparen',
    'bergen',
    'rauben',
    'klingen',
    'ziehen',
    'weichen',
    'dienen',
    'fließen',
    'binden',
    'drücken',
    'rasen',
    'runden',
    'flüchten',
    'liefern',
    'melden',
    'ärgern',
    'sichern',
    'töten',
    'sterben',
    'freuen',
    'süßen',
    'planen',
    'klettern',
    'streichen',
    'spielen',
    'scheinen',
    'greifen',
    'albern',
    'herzen',
    'mögen',
    'kämpfen',
    'seufzen',
    'trauen',
    'bilden',
    'stammen',
    'führen',
    'wiegen',
    'lügen',
    'tropfen',
    'schmecken',
    'segeln',
    'schrecken',
    'hören',
    'treffen',
    'enden',
    'stören',
    'geben',
    'decken',
    'heulen',
    'schlagen',
    'kichern',
    'schenken',
    'backen',
    'lenken',
    'wahren',
    'kehren',
    'wüten',
    'siegen',
    'boxen',
    'wetten',
    'dringen',
    'fällen',
    'kranken',
    'schämen',
    'wirken',
    'feiern',
    'fahren',
    'tollen',
    'scheiden',
    'zünden',
    'schlucken',
    'pflanzen',
    'bauen',
    'duschen',
    'brauchen',
    'zählen',
    'spüren',
    'stecken',
    'räumen',
    'betteln',
    'regeln',
    'münzen',
    'machen',
    'helfen',
    'parken',
    'platzen',
    'graben',
    'äußern',
    'malen',
    'scheitern',
    'holen',
    'schwören',
    'treiben',
    'formen',
    'wecken',
    'sinken',
    'flehen',
    'sorgen',
    'zeigen',
    'schwächen',
    'wohnen',
    'lockern',
    'meistern',
    'schicken',
    'drehen',
    'gründen',
    'heben',
    'hauen',
    'stoßen',
    'schütteln',
    'locken',
    'lesen',
    'ändern',
    'deuten',
    'wenden',
    'bellen',
    'ruhen',
    'altern',
    'schnellen',
    'ehren',
    'saufen',
    'bitten',
    'klappen',
    'gleichen',
    'wachsen',
    'erben',
    'schwingen',
    'feuern',
    'sprengen',
    'heilen',
    'kümmern',
    'suchen',
    'pfeifen',
    'baden',
    'loben',
    'fangen',
    'stürzen',
    'wundern',
    'werfen',
    'jubeln',
    'rufen',
    'lassen',
    'achten',
    'liegen',
    'folgen',
    'wehtun',
    'leihen',
    'ahnen',
    'werden',
    'knurren',
    'gelten',
    'warnen',
    'testen',
    'reizen',
    'zielen',
    'rechnen',
    'wehen',
    'hupen',
    'trennen',
    'quälen',
    'arten'
)

for ($i = 0; $i -lt $words.Length; $i++) {
    $ws.Cells.Item($i + 2, 2).Value = $words[$i]
}
